$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# New rows of mail-log data to append to the "Logs" sheet.
$newRows = @(
    @{
        Onderwerp  = "Kun jij dit even regelen?"
        Afzender   = "mailmind.test@zohomail.eu"
        Inhoud     = "Testmail #1: Kun jij dit even regelen?"
        Categorie  = "Planning / Afspraak"
        Antwoord   = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
        Tijdstip   = "2025-08-05 17:09:08"
        Beantwoord = "Ja"
        Handmatig  = "Ja"
        Automatisch = "Nee"
        Hergebruikt = "Nee"
    },
    @{
        Onderwerp  = "Kun jij dit even regelen?"
        Afzender   = "mailmind.test@zohomail.eu"
        Inhoud     = "Testmail #1: Kun jij dit even regelen?"
        Categorie  = "Planning / Afspraak"
        Antwoord   = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
        Tijdstip   = "2025-08-05 17:09:10"
        Beantwoord = "Ja"
        Handmatig  = "Ja"
        Automatisch = "Nee"
        Hergebruikt = "Nee"
    }
)

$startRow = $logs.Cells.Item($logs.Rows.Count, 1).End(-4162).Row + 1
if ($startRow -lt 2) { $startRow = 2 }

$r = $startRow
foreach ($row in $newRows) {
    $logs.Cells.Item($r, 1).Value = $row.Onderwerp
    $logs.Cells.Item($r, 2).Value = $row.Afzender
    $logs.Cells.Item($r, 3).Value = $row.Inhoud
    $logs.Cells.Item($r, 4).Value = $row.Categorie
    $logs.Cells.Item($r, 5).Value = $row.Antwoord
    $logs.Cells.Item($r, 6).Value = $row.Tijdstip
    $logs.Cells.Item($r, 7).Value = $row.Beantwoord
    $logs.Cells.Item($r, 8).Value = $row.Handmatig
    $logs.Cells.Item($r, 9).Value = $row.Automatisch
    $logs.Cells.Item($r, 10).Value = $row.Hergebruikt
    $r = $r + 1
}

# Update the Dashboard summary count for "Planning / Afspraak" (6 -> 8).
$dashboard.Cells.Item(2, 2).Value = 8

# Extend the conditional-formatting ranges (D/G/H/I/J) to cover the new rows.
$lastRow = $r - 1
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "12")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "$lastRow")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}
